$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 401.875
$ws.Range("I2").Value = 392
$ws.Range("J2").Value = 423.6
$ws.Range("K2").Value = 392
$ws.Range("L2").Value = 423.6
$ws.Range("M2").Value = -279
$ws.Range("N2").Value = -649.6
$ws.Range("H33").Value = 492.22223
$ws.Range("I33").Value = 398.85715
$ws.Range("K33").Value = 398.85715
$ws.Range("M33").Value = -169.85715
$ws.Range("H53").Value = 2927
$ws.Range("J53").Value = 4416.3
$ws.Range("L53").Value = 4416.3
$ws.Range("N53").Value = -5690.3
$ws.Range("H55").Value = 1485.091
$ws.Range("I55").Value = 426.2
$ws.Range("J55").Value = 3754.1428
$ws.Range("K55").Value = 426.2
$ws.Range("L55").Value = 3754.1428
$ws.Range("M55").Value = -212.2
$ws.Range("N55").Value = -4182.1428
$ws.Range("H76").Value = 5559.0713
$ws.Range("I76").Value = 4369.4116
$ws.Range("J76").Value = 7397.636
$ws.Range("K76").Value = 4369.4116
$ws.Range("L76").Value = 7397.636
$ws.Range("M76").Value = -4054.4116
$ws.Range("N76").Value = -8027.636
$ws.Range("H79").Value = 5559.0713
$ws.Range("I79").Value = 4369.4116
$ws.Range("J79").Value = 7397.636
$ws.Range("K79").Value = 4369.4116
$ws.Range("L79").Value = 7397.636
$ws.Range("M79").Value = -3277.4116
$ws.Range("N79").Value = -9581.636
$ws.Range("H112").Value = 1304.8462
$ws.Range("J112").Value = 1356.6818
$ws.Range("L112").Value = 4070.0454
$ws.Range("N112").Value = -6286.0454
$ws.Range("H125").Value = 2936.6
$ws.Range("J125").Value = 2936.6
$ws.Range("L125").Value = 26429.4
$ws.Range("N125").Value = -31349.4
$ws.Range("H132").Value = 2500.8333
$ws.Range("I132").Value = 2201
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 6603
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -4073
$ws.Range("N132").Value = -17060
$ws.Range("H137").Value = 3184.7856
$ws.Range("I137").Value = 2327.4443
$ws.Range("J137").Value = 3348.9575
$ws.Range("K137").Value = 6982.3329
$ws.Range("L137").Value = 10046.8725
$ws.Range("M137").Value = -4432.3329
$ws.Range("N137").Value = -15146.8725

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3665.9019
$ws.Range("I32").Value = 2540.8542
$ws.Range("K32").Value = 2540.8542
$ws.Range("M32").Value = -2253.8542
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("M37").ClearContents()
$ws.Range("H43").Value = 50032
$ws.Range("J43").Value = 56665.332
$ws.Range("L43").Value = 56665.332
$ws.Range("N43").Value = -57291.332
$ws.Range("H44").Value = 15000
$ws.Range("I44").Value = 15000
$ws.Range("K44").Value = 15000
$ws.Range("M44").Value = -14512
$ws.Range("H45").Value = 71430730
$ws.Range("I45").Value = 71430730
$ws.Range("K45").Value = 71430730
$ws.Range("M45").Value = -71430353
$ws.Range("H61").Value = 5591.8335
$ws.Range("I61").Value = 4838.0356
$ws.Range("K61").Value = 4838.0356
$ws.Range("M61").Value = -4626.0356
$ws.Range("H80").Value = 104799
$ws.Range("I80").Value = 77000
$ws.Range("J80").Value = 123331.664
$ws.Range("K80").Value = 77000
$ws.Range("L80").Value = 123331.664
$ws.Range("M80").Value = -76002
$ws.Range("N80").Value = -125327.664
$ws.Range("H83").Value = 104799
$ws.Range("I83").Value = 77000
$ws.Range("J83").Value = 123331.664
$ws.Range("K83").Value = 231000
$ws.Range("L83").Value = 369994.992
$ws.Range("M83").Value = -226008
$ws.Range("N83").Value = -379978.992
$ws.Range("H136").Value = 5591.8335
$ws.Range("I136").Value = 4838.0356
$ws.Range("K136").Value = 14514.1068
$ws.Range("M136").Value = -11964.1068
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").ClearContents()
$ws.Range("N139").Value = 0

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3859.2173
$ws.Range("I20").Value = 3834
$ws.Range("J20").Value = 3916.8572
$ws.Range("K20").Value = 3834
$ws.Range("L20").Value = 3916.8572
$ws.Range("M20").Value = -3587
$ws.Range("N20").Value = -4410.8572
$ws.Range("H22").Value = 866.3333
$ws.Range("J22").Value = 1000
$ws.Range("L22").Value = 1000
$ws.Range("N22").Value = -1346
$ws.Range("H94").Value = 2738.75
$ws.Range("I94").Value = 2660.4546
$ws.Range("J94").Value = 2911
$ws.Range("K94").Value = 2660.4546
$ws.Range("L94").Value = 2911
$ws.Range("M94").Value = -2209.4546
$ws.Range("N94").Value = -3813
$ws.Range("H105").Value = 12610.458
$ws.Range("I105").Value = 12785.529
$ws.Range("J105").Value = 12185.286
$ws.Range("K105").Value = 12785.529
$ws.Range("L105").Value = 12185.286
$ws.Range("M105").Value = -11038.529
$ws.Range("N105").Value = -15679.286
$ws.Range("H134").Value = 3596.5
$ws.Range("I134").Value = 1579.2307
$ws.Range("K134").Value = 4737.6921
$ws.Range("M134").Value = -2202.6921

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 5986.5
$ws.Range("I86").Value = 6270.2856
$ws.Range("K86").Value = 6270.2856
$ws.Range("M86").Value = -5147.2856
$ws.Range("H89").Value = 5986.5
$ws.Range("I89").Value = 6270.2856
$ws.Range("K89").Value = 31351.428
$ws.Range("M89").Value = -25735.428

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 283.1
$ws.Range("I7").Value = 311.57144
$ws.Range("J7").Value = 216.66667
$ws.Range("K7").Value = 934.71432
$ws.Range("L7").Value = 650.00001
$ws.Range("M7").Value = -822.71432
$ws.Range("N7").Value = -874.00001
$ws.Range("H92").Value = 4645.5
$ws.Range("I92").Value = 2909
$ws.Range("J92").Value = 7250.25
$ws.Range("K92").Value = 8727
$ws.Range("L92").Value = 21750.75
$ws.Range("M92").Value = -7479
$ws.Range("N92").Value = -24246.75
$ws.Range("H132").Value = 4015.4285
$ws.Range("I132").Value = 2840.4614
$ws.Range("K132").Value = 25564.1526
$ws.Range("M132").Value = -23034.1526

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 720.9375
$ws.Range("I2").Value = 201.6
$ws.Range("K2").Value = 201.6
$ws.Range("M2").Value = -88.59999999999999
$ws.Range("H70").Value = 14222.385
$ws.Range("I70").Value = 13413.571
$ws.Range("J70").Value = 15166
$ws.Range("K70").Value = 13413.571
$ws.Range("L70").Value = 15166
$ws.Range("M70").Value = -13143.571
$ws.Range("N70").Value = -15706
$ws.Range("H73").Value = 14222.385
$ws.Range("I73").Value = 13413.571
$ws.Range("J73").Value = 15166
$ws.Range("K73").Value = 13413.571
$ws.Range("L73").Value = 15166
$ws.Range("M73").Value = -12477.571
$ws.Range("N73").Value = -17038
$ws.Range("H122").Value = 5854.7144
$ws.Range("I122").Value = 5136.909
$ws.Range("J122").Value = 6644.3
$ws.Range("K122").Value = 15410.727
$ws.Range("L122").Value = 19932.9
$ws.Range("M122").Value = -12960.727
$ws.Range("N122").Value = -24832.9
$ws.Range("H126").Value = 3871.476
$ws.Range("I126").Value = 2248.8
$ws.Range("J126").Value = 5346.636
$ws.Range("K126").Value = 6746.400000000001
$ws.Range("L126").Value = 16039.908
$ws.Range("M126").Value = -4276.400000000001
$ws.Range("N126").Value = -20979.908
$ws.Range("H132").Value = 3472.4285
$ws.Range("I132").Value = 2725.8333
$ws.Range("K132").Value = 8177.499899999999
$ws.Range("M132").Value = -5647.499899999999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5625.7144
$ws.Range("I7").Value = 4255.483
$ws.Range("K7").Value = 4255.483
$ws.Range("M7").Value = -4143.483
$ws.Range("H82").Value = 5570.8335
$ws.Range("I82").Value = 1293.8889
$ws.Range("K82").Value = 1293.8889
$ws.Range("M82").Value = -932.8888999999999
$ws.Range("H85").Value = 5570.8335
$ws.Range("I85").Value = 1293.8889
$ws.Range("K85").Value = 1293.8889
$ws.Range("M85").Value = -45.88889999999992
$ws.Range("H100").Value = 13374.6875
$ws.Range("I100").Value = 8326.666999999999
$ws.Range("K100").Value = 8326.666999999999
$ws.Range("M100").Value = -7785.666999999999
$ws.Range("H126").Value = 5625.7144
$ws.Range("I126").Value = 4255.483
$ws.Range("K126").Value = 12766.449
$ws.Range("M126").Value = -10296.449
$ws.Range("H136").Value = 3045.7942
$ws.Range("I136").Value = 1767.4783
$ws.Range("J136").Value = 5718.636
$ws.Range("K136").Value = 5302.4349
$ws.Range("L136").Value = 17155.908
$ws.Range("M136").Value = -2752.4349
$ws.Range("N136").Value = -22255.908

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 27173334
$ws.Range("J5").Value = 32606000
$ws.Range("L5").Value = 32606000
$ws.Range("N5").Value = -32606224
$ws.Range("H126").Value = 1537.5625
$ws.Range("I126").Value = 1261.8889
$ws.Range("J126").Value = 1892
$ws.Range("K126").Value = 3785.6667
$ws.Range("L126").Value = 5676
$ws.Range("M126").Value = -1315.6667
$ws.Range("N126").Value = -10616
$ws.Range("H136").Value = 4144.44
$ws.Range("I136").Value = 3148.9412
$ws.Range("J136").Value = 6259.875
$ws.Range("K136").Value = 9446.8236
$ws.Range("L136").Value = 18779.625
$ws.Range("M136").Value = -6896.8236
$ws.Range("N136").Value = -23879.625
